# Weekly update: a new week of "Ají" price data (variety "Inferno", Primera
# and Segunda quality) is prepended right after the existing row for
# 2021-08-16 (row 11). All the historical rows that used to start at row 12
# shift down by two rows (old row 12 -> new row 14, ... old row 100 -> new
# row 102).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 12:13, pushing every row from 12 downward to 14
# onward. This mirrors Excel's own "Insert Copied/Blank Rows" behaviour,
# carrying the date-column (D) number format down into the new rows and
# automatically growing the sheet's used range / <dimension> to A1:R102.
$ws.Range("A12:A13").EntireRow.Insert()

# --- New row 12: Inferno / Primera, week of 2021-10-04 ---
$ws.Cells.Item(12, 1).Value = 8
$ws.Cells.Item(12, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(12, 3).Value = "Coquimbo"
$ws.Cells.Item(12, 4).Value = 44473
$ws.Cells.Item(12, 5).Value = 4
$ws.Cells.Item(12, 6).Value = 100112021
$ws.Cells.Item(12, 7).Value = "Ají"
$ws.Cells.Item(12, 8).Value = "Inferno"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 500
$ws.Cells.Item(12, 11).Value = 39000
$ws.Cells.Item(12, 12).Value = 40000
$ws.Cells.Item(12, 13).Value = 39500
$ws.Cells.Item(12, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(12, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(12, 16).Value = 3292
$ws.Cells.Item(12, 17).Value = 12
$ws.Cells.Item(12, 18).Value = "Hortaliza"

# --- New row 13: Inferno / Segunda, week of 2021-10-04 ---
$ws.Cells.Item(13, 1).Value = 8
$ws.Cells.Item(13, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(13, 3).Value = "Coquimbo"
$ws.Cells.Item(13, 4).Value = 44473
$ws.Cells.Item(13, 5).Value = 4
$ws.Cells.Item(13, 6).Value = 100112021
$ws.Cells.Item(13, 7).Value = "Ají"
$ws.Cells.Item(13, 8).Value = "Inferno"
$ws.Cells.Item(13, 9).Value = "Segunda"
$ws.Cells.Item(13, 10).Value = 360
$ws.Cells.Item(13, 11).Value = 33000
$ws.Cells.Item(13, 12).Value = 34000
$ws.Cells.Item(13, 13).Value = 33500
$ws.Cells.Item(13, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(13, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(13, 16).Value = 2792
$ws.Cells.Item(13, 17).Value = 12
$ws.Cells.Item(13, 18).Value = "Hortaliza"

Write-Output "done"
